$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.42884377887058
$ws.Range("C2").Value = 8.733351859699688
$ws.Range("D2").Value = 5.99240110028064
$ws.Range("E2").Value = 11.88201854985191
$ws.Range("G2").Value = 31.90462562968723
$ws.Range("H2").Value = 15.15969737437433
$ws.Range("I2").Value = 23.82754678969025
$ws.Range("K2").Value = 8.734072931944917
$ws.Range("L2").Value = 9.693446550155713
$ws.Range("O2").Value = 23.4799202700349

$ws.Range("B3").Value = 11.12446420635811
$ws.Range("C3").Value = 8.714002696253743
$ws.Range("D3").Value = 5.874675357440139
$ws.Range("E3").Value = 11.89674192154968
$ws.Range("G3").Value = 31.9973920134633
$ws.Range("H3").Value = 15.20993769413628
$ws.Range("I3").Value = 23.92563600816627
$ws.Range("K3").Value = 8.511829824201762
$ws.Range("L3").Value = 9.677160502835264
$ws.Range("O3").Value = 23.56268646965013

$ws.Range("B4").Value = 10.93476586125775
$ws.Range("C4").Value = 8.702231951296076
$ws.Range("D4").Value = 5.802917148986147
$ws.Range("E4").Value = 11.90798196259104
$ws.Range("G4").Value = 32.06398089040322
$ws.Range("H4").Value = 15.24314300811064
$ws.Range("I4").Value = 23.99032203942776
$ws.Range("K4").Value = 8.373280123315514
$ws.Range("L4").Value = 9.668806924339659
$ws.Range("O4").Value = 23.61835063199169

$ws.Range("B5").Value = 10.85687026102063
$ws.Range("C5").Value = 8.69746432246297
$ws.Range("D5").Value = 5.773852745284229
$ws.Range("E5").Value = 11.91311558509576
$ws.Range("G5").Value = 32.09352844890649
$ws.Range("H5").Value = 15.25726734194343
$ws.Range("I5").Value = 24.01780275552591
$ws.Range("K5").Value = 8.316373492351275
$ws.Range("L5").Value = 9.66581947793771
$ws.Range("O5").Value = 23.64225068795184

$ws.Range("B6").Value = 10.84390346129705
$ws.Range("C6").Value = 8.696674457734074
$ws.Range("D6").Value = 5.769038688260827
$ws.Range("E6").Value = 11.91400143065637
$ws.Range("G6").Value = 32.09858018952769
$ws.Range("H6").Value = 15.25964848869127
$ws.Range("I6").Value = 24.02243357442448
$ws.Range("K6").Value = 8.306899625248883
$ws.Range("L6").Value = 9.665348658034475
$ws.Range("O6").Value = 23.64629268704587

$ws.Range("B7").Value = 10.93371757538925
$ws.Range("C7").Value = 8.702167533873995
$ws.Range("D7").Value = 5.802524396641746
$ws.Range("E7").Value = 11.90804895650936
$ws.Range("G7").Value = 32.06436962587495
$ws.Range("H7").Value = 15.24333109309625
$ws.Range("I7").Value = 23.99068811733377
$ws.Range("K7").Value = 8.372514360424624
$ws.Range("L7").Value = 9.66876494391936
$ws.Range("O7").Value = 23.61866803373964

$ws.Range("B8").Value = 11.324542968533
$ws.Range("C8").Value = 8.726658371416505
$ws.Range("D8").Value = 5.951727198732017
$ws.Range("E8").Value = 11.88663866260831
$ws.Range("G8").Value = 31.93460729869559
$ws.Range("H8").Value = 15.17653097073952
$ws.Range("I8").Value = 23.86044216610481
$ws.Range("K8").Value = 8.657923932127169
$ws.Range("L8").Value = 9.687491217485599
$ws.Range("O8").Value = 23.50745123970087

$ws.Range("B9").Value = 12.06399541578113
$ws.Range("C9").Value = 8.775492149426704
$ws.Range("D9").Value = 6.246525281230951
$ws.Range("E9").Value = 11.86210102307946
$ws.Range("G9").Value = 31.75694520777366
$ws.Range("H9").Value = 15.06423891659853
$ws.Range("I9").Value = 23.6404351897891
$ws.Range("K9").Value = 9.197750548144432
$ws.Range("L9").Value = 9.737141040961877
$ws.Range("O9").Value = 23.32788682364752

$ws.Range("B10").Value = 12.58508331737763
$ws.Range("C10").Value = 8.811778655675537
$ws.Range("D10").Value = 6.461914830710781
$ws.Range("E10").Value = 11.85469374684624
$ws.Range("G10").Value = 31.67372503537339
$ws.Range("H10").Value = 14.99313488709337
$ws.Range("I10").Value = 23.50041524557287
$ws.Range("K10").Value = 9.578215768643059
$ws.Range("L10").Value = 9.781303959639475
$ws.Range("O10").Value = 23.21956464485018

$ws.Range("B11").Value = 12.81625284130904
$ws.Range("C11").Value = 8.828356850430625
$ws.Range("D11").Value = 6.559132611011775
$ws.Range("E11").Value = 11.8536235561098
$ws.Range("G11").Value = 31.64622544009425
$ws.Range("H11").Value = 14.96326169561904
$ws.Range("I11").Value = 23.4414190920533
$ws.Range("K11").Value = 9.747038444089101
$ws.Range("L11").Value = 9.803016008395263
$ws.Range("O11").Value = 23.17543470835703

$ws.Range("B12").Value = 12.90286670891712
$ws.Range("C12").Value = 8.83464321903562
$ws.Range("D12").Value = 6.595796903690247
$ws.Range("E12").Value = 11.85354809671238
$ws.Range("G12").Value = 31.63730709744415
$ws.Range("H12").Value = 14.95230501324753
$ws.Range("I12").Value = 23.41975540880477
$ws.Range("K12").Value = 9.810299464797463
$ws.Range("L12").Value = 9.81146656597301
$ws.Range("O12").Value = 23.15946573844872

$ws.Range("B13").Value = 12.88425520050972
$ws.Range("C13").Value = 8.833288983094327
$ws.Range("D13").Value = 6.587907877220328
$ws.Range("E13").Value = 11.85354969579615
$ws.Range("G13").Value = 31.6391612437932
$ws.Range("H13").Value = 14.95464891139157
$ws.Range("I13").Value = 23.42439094877105
$ws.Range("K13").Value = 9.796705655755327
$ws.Range("L13").Value = 9.809636492559548
$ws.Range("O13").Value = 23.16287190607675

$ws.Range("B14").Value = 12.82339758768561
$ws.Range("C14").Value = 8.828873871787529
$ws.Range("D14").Value = 6.562152231918914
$ws.Range("E14").Value = 11.85361074386819
$ws.Range("G14").Value = 31.64546172816234
$ws.Range("H14").Value = 14.96235315437101
$ws.Range("I14").Value = 23.43962323585379
$ws.Range("K14").Value = 9.752256671709928
$ws.Range("L14").Value = 9.803706685037215
$ws.Range("O14").Value = 23.17410604783247

$ws.Range("B15").Value = 12.78599782681235
$ws.Range("C15").Value = 8.826170553914075
$ws.Range("D15").Value = 6.546355472307014
$ws.Range("E15").Value = 11.85369105880305
$ws.Range("G15").Value = 31.64951582373572
$ws.Range("H15").Value = 14.96711854500772
$ws.Range("I15").Value = 23.4490416336553
$ws.Range("K15").Value = 9.724941727520957
$ws.Range("L15").Value = 9.800104146529973
$ws.Range("O15").Value = 23.18108398377112

$ws.Range("B16").Value = 12.56985059536287
$ws.Range("C16").Value = 8.810696522446419
$ws.Range("D16").Value = 6.455542567757025
$ws.Range("E16").Value = 11.85480988730996
$ws.Range("G16").Value = 31.67573113432312
$ws.Range("H16").Value = 14.99513693140704
$ws.Range("I16").Value = 23.5043654431562
$ws.Range("K16").Value = 9.567092253521247
$ws.Range("L16").Value = 9.77991728171108
$ws.Range("O16").Value = 23.22255237408346

$ws.Range("B17").Value = 12.43568700364621
$ws.Range("C17").Value = 8.801220769368191
$ws.Range("D17").Value = 6.399606883931186
$ws.Range("E17").Value = 11.85608463417444
$ws.Range("G17").Value = 31.69447062001701
$ws.Range("H17").Value = 15.01295862654233
$ws.Range("I17").Value = 23.5395092554553
$ws.Range("K17").Value = 9.469125386001606
$ws.Range("L17").Value = 9.767945581872089
$ws.Range("O17").Value = 23.24931147162657

$ws.Range("B18").Value = 12.35797089390208
$ws.Range("C18").Value = 8.795777377995481
$ws.Range("D18").Value = 6.36736451098048
$ws.Range("E18").Value = 11.85703432874634
$ws.Range("G18").Value = 31.70622370836976
$ws.Range("H18").Value = 15.02344189963694
$ws.Range("I18").Value = 23.56016533134882
$ws.Range("K18").Value = 9.412380218488712
$ws.Range("L18").Value = 9.761212749524677
$ws.Range("O18").Value = 23.26518688009917

$ws.Range("B19").Value = 12.33156586366365
$ws.Range("C19").Value = 8.793935560063375
$ws.Range("D19").Value = 6.356437093445042
$ws.Range("E19").Value = 11.85739308528463
$ws.Range("G19").Value = 31.71037030529711
$ws.Range("H19").Value = 15.02703131955397
$ws.Range("I19").Value = 23.56723505400544
$ws.Range("K19").Value = 9.393100872827228
$ws.Range("L19").Value = 9.758959536615929
$ws.Range("O19").Value = 23.27064513336639

$ws.Range("B20").Value = 12.45002634471203
$ws.Range("C20").Value = 8.802228786675787
$ws.Range("D20").Value = 6.405568838699736
$ws.Range("E20").Value = 11.85592653454246
$ws.Range("G20").Value = 31.69237485434918
$ws.Range("H20").Value = 15.01103739166416
$ws.Range("I20").Value = 23.53572235408373
$ws.Range("K20").Value = 9.479595664638856
$ws.Range("L20").Value = 9.769204189998087
$ws.Range("O20").Value = 23.24641278347527

$ws.Range("B21").Value = 12.84129865034672
$ws.Range("C21").Value = 8.83017047802635
$ws.Range("D21").Value = 6.569721668072096
$ws.Range("E21").Value = 11.85358386959069
$ws.Range("G21").Value = 31.64357050423299
$ws.Range("H21").Value = 14.96008057832298
$ws.Range("I21").Value = 23.4351307648662
$ws.Range("K21").Value = 9.765330986680107
$ws.Range("L21").Value = 9.805442244883029
$ws.Range("O21").Value = 23.17078615242011

$ws.Range("B22").Value = 13.09159397291834
$ws.Range("C22").Value = 8.848481125348526
$ws.Range("D22").Value = 6.676113872686845
$ws.Range("E22").Value = 11.85397460360567
$ws.Range("G22").Value = 31.62039091610522
$ws.Range("H22").Value = 14.92885039551974
$ws.Range("I22").Value = 23.37333423557942
$ws.Range("K22").Value = 9.948156284482073
$ws.Range("L22").Value = 9.830456496510887
$ws.Range("O22").Value = 23.12568607037382

$ws.Range("B23").Value = 12.95852750289815
$ws.Range("C23").Value = 8.838704452853245
$ws.Range("D23").Value = 6.619424563499216
$ws.Range("E23").Value = 11.853590543256
$ws.Range("G23").Value = 31.63196306614978
$ws.Range("H23").Value = 14.94532880788434
$ws.Range("I23").Value = 23.40595479475497
$ws.Range("K23").Value = 9.850955124446568
$ws.Range("L23").Value = 9.816985752603751
$ws.Range("O23").Value = 23.14936032472789

$ws.Range("B24").Value = 12.44354534047298
$ws.Range("C24").Value = 8.801773048392187
$ws.Range("D24").Value = 6.402873698160183
$ws.Range("E24").Value = 11.8559973359416
$ws.Range("G24").Value = 31.69331929934688
$ws.Range("H24").Value = 15.01190524264421
$ws.Range("I24").Value = 23.53743300608879
$ws.Range("K24").Value = 9.474863363941672
$ws.Range("L24").Value = 9.768634706135387
$ws.Range("O24").Value = 23.24772175033089

$ws.Range("B25").Value = 11.86746259010078
$ws.Range("C25").Value = 8.762203460404566
$ws.Range("D25").Value = 6.166806411142129
$ws.Range("E25").Value = 11.86687223082161
$ws.Range("G25").Value = 31.79673169158723
$ws.Range("H25").Value = 15.09261514502806
$ws.Range("I25").Value = 23.69615900570065
$ws.Range("K25").Value = 9.054275421007532
$ws.Range("L25").Value = 9.722344519634307
$ws.Range("O25").Value = 23.37232668521125
